# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted as a new row 15 (pushing the
# existing rows 15-44 down to 16-45, which Excel's native Rows(...).Insert()
# does automatically, preserving all of their original values/formatting).
# The freshly inserted row 15 is then populated with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; rows 15-44 shift down to 16-45.
$ws.Rows(15).Insert()

# Populate the new row 15 with the new observation.
$ws.Cells.Item(15, 1).Value  = 5                                          # A15 Mercado ID
$ws.Cells.Item(15, 2).Value  = "Macroferia Regional de Talca"             # B15 Mercado
$ws.Cells.Item(15, 3).Value  = "Maule"                                    # C15 Región
$ws.Cells.Item(15, 4).Value  = 44414                                      # D15 Fecha
$ws.Cells.Item(15, 5).Value  = 7                                          # E15 Codreg
$ws.Cells.Item(15, 6).Value  = 100112001                                  # F15 Categoría ID
$ws.Cells.Item(15, 7).Value  = "Berenjena"                                # G15 Categoría
$ws.Cells.Item(15, 8).Value  = "Sin especificar"                         # H15 Variedad
$ws.Cells.Item(15, 9).Value  = "Primera"                                  # I15 Calidad
$ws.Cells.Item(15, 10).Value = 300                                        # J15 Volumen
$ws.Cells.Item(15, 11).Value = 7000                                       # K15 Precio mínimo
$ws.Cells.Item(15, 12).Value = 7000                                       # L15 Precio máximo
$ws.Cells.Item(15, 13).Value = 7000                                       # M15 Precio promedio ponderado
$ws.Cells.Item(15, 14).Value = "$/caja 60 unidades"                       # N15 Unidad de comercialización
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"             # O15 Origen
$ws.Cells.Item(15, 16).Value = 117                                        # P15 Precio $/Kg
$ws.Cells.Item(15, 17).Value = 60                                         # Q15 Kg o Unidades
$ws.Cells.Item(15, 18).Value = "Hortaliza"                                # R15 Clasificación
